# Word re-touched every inline picture (marking its run as "no
# proofing", i.e. <w:noProof/> in the run's rPr -- creating the rPr
# if the run didn't already have one) and dropped one stray empty
# paragraph near the end of the body.

$d = $word.ActiveDocument

# 1) Mark every inline picture's range as "no proofing". This inserts
#    <w:noProof/> into the owning run's <w:rPr> (in schema order,
#    right after w:b/w:bCs and before w:color/w:u/...), creating the
#    <w:rPr> wrapper when the run didn't have one yet.
$shapeCount = $d.InlineShapes.Count
for ($i = 1; $i -le $shapeCount; $i++) {
    $d.InlineShapes.Item($i).Range.NoProofing = 1
}

# 2) Remove the extra empty paragraph sandwiched between the two
#    underline-formatted empty paragraphs near the end of the
#    document (identified structurally: an empty, non-underlined
#    paragraph whose immediate neighbours are both empty and
#    underlined -- matches exactly one paragraph in this document).
$paraCount = $d.Paragraphs.Count
for ($i = $paraCount - 1; $i -ge 2; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq [char]13 -and $p.Range.Font.Underline -eq 0) {
        $prev = $d.Paragraphs.Item($i - 1)
        $next = $d.Paragraphs.Item($i + 1)
        if ($prev.Range.Text -eq [char]13 -and $next.Range.Text -eq [char]13 -and `
            $prev.Range.Font.Underline -eq 1 -and $next.Range.Font.Underline -eq 1) {
            $p.Range.Delete()
            break
        }
    }
}
